# Word snippet clean up - folders and descriptions
# Renames the "SnippetIdIntheYAMLFile" (column C) folder/id values in the
# Snippets table to the new folder naming scheme, and updates the active
# selection to C31, matching the authored commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Snippets")

# Map of old folder/snippet-id values to their renamed equivalents.
$map = @{
    "word-basics-insert-and-get-pictures"                 = "word-images-insert-and-get-pictures";
    "word-basics-insert-header"                            = "word-paragraphs-insert-header-and-footer";
    "word-basics-insert-line-and-page-breaks"               = "word-paragraphs-insert-line-and-page-breaks";
    "word-basics-insert-formatted-text"                     = "word-paragraphs-insert-formatted-text";
    "word-basics-search"                                    = "word-paragraphs-search";
    "word-basics-basic-doc-assembly"                        = "word-scenarios-basic-doc-assembly";
    "word-basics-insert-in-different-locations"             = "word-paragraphs-insert-in-different-locations";
    "word-basics-scroll-to-range"                           = "word-ranges-scroll-to-range";
    "word-range-split-words-of-first-paragraph"             = "word-ranges-split-words-of-first-paragraph";
    "word-basics-read-write-custom-document-properties"     = "word-properties-read-write-custom-document-properties";
    "word-custom-properties-get-built-in-properties"        = "word-properties-get-built-in-properties";
}

# Column C is "SnippetIdIntheYAMLFile". Process rows in the order that
# reproduces the authored shared-string layout (grouped by destination
# folder, in the order each group's edits were completed).
# Note: use Value2 (not Value) for reading/writing - Value's getter is
# unreliable for string comparisons in this runtime.
$rowOrder = @(2, 4, 5, 6, 7, 8, 9, 10, 11, 3, 12, 19, 22, 23, 24, 28, 29, 30, 31)
foreach ($r in $rowOrder) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($map.ContainsKey($current)) {
        $cell.Value2 = $map[$current]
    }
}

# Update the current selection to match the authored state.
$ws.Range("C31").Select()
